$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: remove rows that duplicate content now uniquely kept on Sheet2 ---
# (original row numbers, deleted from bottom to top so indices stay valid)
[void]$ws1.Rows.Item(21).Delete()
[void]$ws1.Rows.Item(17).Delete()
[void]$ws1.Rows.Item(16).Delete()
[void]$ws1.Rows.Item(13).Delete()
[void]$ws1.Rows.Item(12).Delete()
[void]$ws1.Rows.Item(11).Delete()
[void]$ws1.Rows.Item(10).Delete()

# --- Sheet2: remove the single row that duplicates content now uniquely kept on Sheet1 ---
[void]$ws2.Rows.Item(15).Delete()

# --- Column C width on Sheet1 widened to match Sheet2's best-fit width ---
$ws1.Columns.Item(3).ColumnWidth = 28.17

# --- Sheet2 gains an explicit page setup (paper size / orientation) ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Selection / active sheet bookkeeping ---
[void]$ws1.Range("A11").Select()
[void]$ws2.Activate()
[void]$ws2.Range("A16").Select()
